$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-CellText "D2" '26.975.98'
Set-CellText "E2" '  +1.18%  '
Set-CellText "D3" '1.846.93'
Set-CellText "E3" '  +1.06%  '
Set-CellText "E4" '  +0.42%  '
Set-CellText "B5" 'USDC'
Set-CellText "C5" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-CellText "D5" '1.012'
Set-CellText "E5" '  +0.40%  '
Set-CellText "B6" 'BNB'
Set-CellText "C6" 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-CellText "D6" '309.99'
Set-CellText "E6" '  +0.04%  '
Set-CellText "D7" '0.4776'
Set-CellText "E7" '  +2.36%  '
Set-CellText "D8" '0.3679'
Set-CellText "E8" '  +2.21%  '
Set-CellText "E9" '  +1.14%  '
Set-CellText "D10" '0.9299'
Set-CellText "E10" '  +2.76%  '
Set-CellText "D11" '19.72'
Set-CellText "E11" '  +1.36%  '
Set-CellText "D12" '0.07719'
Set-CellText "E12" '  +0.22%  '
Set-CellText "D13" '1.806.82'
Set-CellText "E13" '  -1.10%  '
Set-CellText "D14" '5.337'
Set-CellText "E14" '  +1.23%  '
Set-CellText "D15" '6.431'
Set-CellText "E15" '  +0.85%  '
Set-CellText "D16" '88.76'
Set-CellText "E16" '  +1.08%  '
Set-CellText "D17" '1.014'
Set-CellText "E17" '  +0.45%  '
Set-CellText "D18" '0.000008636'
Set-CellText "E18" '  +0.84%  '
Set-CellText "E19" '  +0.36%  '
Set-CellText "D20" '27.016.60'
Set-CellText "E20" '  +1.23%  '
Set-CellText "D21" '14.48'
Set-CellText "E21" '  +1.80%  '
Set-CellText "D22" '5.058'
Set-CellText "E22" '  +0.54%  '
Set-CellText "D23" '10.63'
Set-CellText "E23" '  +0.65%  '
Set-CellText "D24" '1.934'
Set-CellText "E24" '  +1.21%  '
Set-CellText "D25" '152.83'
Set-CellText "E25" '  -0.07%  '
Set-CellText "D26" '18.22'
Set-CellText "E26" '  +1.41%  '
Set-CellText "D27" '2.002'
Set-CellText "E27" '  +0.33%  '
Set-CellText "D28" '114.35'
Set-CellText "E28" '  +0.36%  '
Set-CellText "D29" '4.981'
Set-CellText "E29" '  +2.12%  '
Set-CellText "D30" '0.08898'
Set-CellText "E30" '  +0.89%  '
Set-CellText "D31" '3.306'
Set-CellText "E31" '  +5.80%  '
Set-CellText "D32" '1.176'
Set-CellText "E32" '  +0.80%  '
Set-CellText "D33" '0.7428'
Set-CellText "E33" '  +0.97%  '
Set-CellText "D34" '4.503'
Set-CellText "E34" '  +1.40%  '
Set-CellText "D35" '2.736'
Set-CellText "E35" '  -4.26%  '
Set-CellText "E36" '  +3.34%  '
Set-CellText "D37" '0.01957'
Set-CellText "E37" '  +1.06%  '
Set-CellText "D38" '0.05265'
Set-CellText "E38" '  +1.98%  '
Set-CellText "D39" '2.976'
Set-CellText "E39" '  +1.96%  '
Set-CellText "D40" '0.5216'
Set-CellText "E40" '  +2.91%  '
Set-CellText "D41" '6.999'
Set-CellText "E41" '  +1.74%  '
Set-CellText "E42" '  +0.98%  '
Set-CellText "D43" '8.203'
Set-CellText "E43" '  +1.50%  '
Set-CellText "D44" '10.62'
Set-CellText "E44" '  +5.62%  '
Set-CellText "E45" '  +1.62%  '
Set-CellText "E46" '  +0.44%  '
Set-CellText "D47" '101.97'
Set-CellText "E47" '  +3.58%  '
Set-CellText "D48" '1.609'
Set-CellText "E48" '  +1.98%  '
Set-CellText "D49" '66.03'
Set-CellText "E49" '  +3.29%  '
Set-CellText "E50" '  +0.44%  '
Set-CellText "D51" '0.8857'
Set-CellText "E51" '  +3.76%  '

Write-Host "Crypto list updated."
